$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.1195823126968207
    "D2" = 0.2521840007823073
    "E2" = 0.1817321604552902
    "F2" = 1.002657103835233
    "G2" = 0.4691778255019017
    "H2" = 0.624537289142431
    "I2" = 0.4322752185273107
    "J2" = 0.1767144191186674
    "K2" = 0.836340381242195
    "O2" = 2.145773834768136
    "B3" = 0.1054339025360349
    "D3" = 0.2451045383362214
    "E3" = 0.1771693951812097
    "F3" = 1.005474163352524
    "G3" = 0.4723886240265074
    "H3" = 0.6300606986291726
    "I3" = 0.4391317139909034
    "J3" = 0.1727077392993479
    "K3" = 0.7314491664905063
    "O3" = 2.163876794879457
    "B4" = 0.09673000671905641
    "D4" = 0.2408494425518768
    "E4" = 0.1744545561820914
    "F4" = 1.007840017497806
    "G4" = 0.4747494121901994
    "H4" = 0.6337673373468036
    "I4" = 0.4436260215305605
    "J4" = 0.1703493160176635
    "K4" = 0.6667970383132911
    "O4" = 2.1764681691173
    "B5" = 0.09317916990835329
    "D5" = 0.2391386688406669
    "E5" = 0.1733701089234394
    "F5" = 1.008964067103889
    "G5" = 0.475809201007614
    "H5" = 0.6353571048546343
    "I5" = 0.4455290126485583
    "J5" = 0.1694138346030272
    "K5" = 0.6403902091896612
    "O5" = 2.181970118322823
    "B6" = 0.09258932698703859
    "D6" = 0.2388560024955098
    "E6" = 0.1731913602184925
    "F6" = 1.00916037423243
    "G6" = 0.4759910773332479
    "H6" = 0.6356258733323159
    "I6" = 0.4458493231257474
    "J6" = 0.1692600452054265
    "K6" = 0.6360017743506319
    "O6" = 2.182906102811586
    "B7" = 0.09668213443382001
    "D7" = 0.2408262762763087
    "E7" = 0.174439842301787
    "F7" = 1.007854529247105
    "G7" = 0.4747633092997532
    "H7" = 0.6337884564822076
    "I7" = 0.4436513963183106
    "J7" = 0.1703365961229437
    "K7" = 0.666441149236249
    "O7" = 2.176540869137625
    "B8" = 0.114707569944656
    "D8" = 0.2497240511433176
    "E8" = 0.1801409585482645
    "F8" = 1.003496372176514
    "G8" = 0.470204006296548
    "H8" = 0.6263763390845583
    "I8" = 0.4345803200484681
    "J8" = 0.1753118216176262
    "K8" = 0.8002265870659357
    "O8" = 2.151709185957401
    "B9" = 0.1499124367479823
    "D9" = 0.2678949077518951
    "E9" = 0.1920067519180222
    "F9" = 1.000000015184035
    "G9" = 0.464359894327302
    "H9" = 0.6143420847684027
    "I9" = 0.4190478570801979
    "J9" = 0.185874824323875
    "K9" = 1.06053556122373
    "O9" = 2.114741449861313
    "B10" = 0.1756789511267129
    "D10" = 0.2816796712387202
    "E10" = 0.2011409425320565
    "F10" = 1.000514636430772
    "G10" = 0.4619648098734217
    "H10" = 0.6070246872467777
    "I10" = 0.4090103917808978
    "J10" = 0.1941278648999969
    "K10" = 1.250459902285456
    "O10" = 2.094752739223196
    "B11" = 0.1873771826285235
    "D11" = 0.2880439751445749
    "E11" = 0.2053864299151087
    "F11" = 1.001419410855704
    "G11" = 0.4612898219435095
    "H11" = 0.6040267002922377
    "I11" = 0.4047422882102474
    "J11" = 0.1979895144136634
    "K11" = 1.336558304482367
    "O11" = 2.087221467412178
    "B12" = 0.1918034373503872
    "D12" = 0.2904672945047366
    "E12" = 0.2070070191796773
    "F12" = 1.001858527709885
    "G12" = 0.4610940116827891
    "H12" = 0.6029389963966167
    "I12" = 0.4031689136304877
    "J12" = 0.1994672427019282
    "K12" = 1.369116820938189
    "O12" = 2.084594518215368
    "B13" = 0.1908503294635295
    "D13" = 0.2899448004841645
    "E13" = 0.2066574229815643
    "F13" = 1.001759663356964
    "G13" = 0.4611335207399208
    "H13" = 0.6031711371087454
    "I13" = 0.4035058614120377
    "J13" = 0.1991483027050123
    "K13" = 1.362106799581795
    "O13" = 2.085150266178459
    "B14" = 0.1877414073186401
    "D14" = 0.2882430775185583
    "E14" = 0.2055194983255006
    "F14" = 1.001453602971523
    "G14" = 0.4612725131397326
    "H14" = 0.6039362608354537
    "I14" = 0.4046119866087921
    "J14" = 0.1981107794015031
    "K14" = 1.339237828486148
    "O14" = 2.087000835446645
    "B15" = 0.1858366237501343
    "D15" = 0.2872024499087331
    "E15" = 0.2048241665084092
    "F15" = 1.001278700750831
    "G15" = 0.4613654421774953
    "H15" = 0.6044111163433996
    "I15" = 0.4052951027373837
    "J15" = 0.197477271936549
    "K15" = 1.325223991520716
    "O15" = 2.088163673722562
    "B16" = 0.174913954653789
    "D16" = 0.2812656180415019
    "E16" = 0.2008653010486228
    "F16" = 1.000469007790713
    "G16" = 0.4620172792366759
    "H16" = 0.6072272687737978
    "I16" = 0.4092953198020357
    "J16" = 0.1938776541152691
    "K16" = 1.244826973146246
    "O16" = 2.095276385957703
    "B17" = 0.1682071300009795
    "D17" = 0.2776474121961172
    "E17" = 0.1984597435401412
    "F17" = 1.000144103102556
    "G17" = 0.4625234654489105
    "H17" = 0.6090396027686751
    "I17" = 0.4118256547876271
    "J17" = 0.191696870207835
    "K17" = 1.19542791879644
    "O17" = 2.100040084561243
    "B18" = 0.1643473863267246
    "D18" = 0.2755751302083809
    "E18" = 0.1970846339558108
    "F18" = 1.000020349909484
    "G18" = 0.4628536188253562
    "H18" = 0.6101131390119292
    "I18" = 0.4133090839273734
    "J18" = 0.1904526420393893
    "K18" = 1.166986856237884
    "O18" = 2.102927007578245
    "B19" = 0.1630401821499134
    "D19" = 0.2748750097602652
    "E19" = 0.1966205080686194
    "F19" = 0.999989288954886
    "G19" = 0.4629720972728819
    "H19" = 0.6104819661970708
    "I19" = 0.4138161643237268
    "J19" = 0.1900331033554892
    "K19" = 1.157352449523898
    "O19" = 2.103929699722869
    "B20" = 0.168921308144121
    "D20" = 0.2780316652679886
    "E20" = 0.1987149396690597
    "F20" = 1.00017215631064
    "G20" = 0.4624655425137263
    "H20" = 0.6088434548475306
    "I20" = 0.4115533934181546
    "J20" = 0.1919279731306176
    "K20" = 1.200689448767719
    "O20" = 2.099517766927903
    "B21" = 0.1886546734540673
    "D21" = 0.2887425552102343
    "E21" = 0.2058533844171961
    "F21" = 1.001540880941249
    "G21" = 0.4612300634974389
    "H21" = 0.6037102343735796
    "I21" = 0.4042859274120953
    "J21" = 0.1984151072417859
    "K21" = 1.345956239509349
    "O21" = 2.086451169245805
    "B22" = 0.2015303824476007
    "D22" = 0.2958201653695198
    "E22" = 0.2105940187323014
    "F22" = 1.002997913539531
    "G22" = 0.4607712092570182
    "H22" = 0.6006326492681922
    "I22" = 0.3997860667377076
    "J22" = 0.2027446074681762
    "K22" = 1.440632773754317
    "O22" = 2.07922293616781
    "B23" = 0.1946604088543324
    "D23" = 0.2920356788961556
    "E23" = 0.2080569907984753
    "F23" = 1.002168781387837
    "G23" = 0.4609841516168984
    "H23" = 0.6022498408653689
    "I23" = 0.402164861423099
    "J23" = 0.2004256646626601
    "K23" = 1.390126934831983
    "O23" = 2.082960639021593
    "B24" = 0.1685984403350176
    "D24" = 0.2778579199206206
    "E24" = 0.1985995409572965
    "F24" = 1.000159277098028
    "G24" = 0.4624916075674292
    "H24" = 0.6089320348867773
    "I24" = 0.4116763933687331
    "J24" = 0.1918234618024997
    "K24" = 1.198310838273812
    "O24" = 2.099753445097605
    "B25" = 0.1404050716798224
    "D25" = 0.2629024695460203
    "E25" = 0.1887234892853797
    "F25" = 1.000404737651415
    "G25" = 0.4656082438106921
    "H25" = 0.6173299779809014
    "I25" = 0.4230083994693743
    "J25" = 0.1829308286121147
    "K25" = 0.9903420687317066
    "O25" = 2.123484417053447
}

foreach ($cellRef in $values.Keys) {
    $ws.Range($cellRef).Value = $values[$cellRef]
}
